$wb = $excel.ActiveWorkbook

# --- SignUpTest: update a few "Expected code" cells ---
$ws1 = $wb.Worksheets.Item("SignUpTest")
$ws1.Range("J5").Value = "201"
$ws1.Range("J20").Value = "500"
$ws1.Range("J21").Value = "500"

# --- UpdateUserTest: remove the "WithId" scenario row (row 2) ---
$ws2 = $wb.Worksheets.Item("UpdateUserTest")
$ws2.Rows.Item(2).Delete()

# --- restore cursor/selection positions ---
$ws1.Activate()
$ws1.Range("J21").Select()

$ws2.Activate()
$ws2.Range("L20").Select()
